$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets - strip spaces / punctuation from the tab names
# ---------------------------------------------------------------------------
$renames = @{
    "802.1x Security"          = "802xSecurity"
    "Computer Name"             = "ComputerName"
    "Ethernet Setup"            = "EthernetSetup"
    "Wireless Properties"       = "WirelessProperties"
    "Wireless Setup"            = "WirelessSetup"
    "Display Settings"          = "DisplaySettings"
    "Keyboard Settings"         = "KeyboardSettings"
    "Mouse Settings"            = "MouseSettings"
    "Power Option"               = "PowerOption"
    "Add Printer"                = "AddPrinter"
    "Date & Time"                = "DateAndTime"
    "Region & Location"          = "RegionAndLocation"
    "Application Command"        = "ApplicationCommand"
    "Environment Variable"       = "EnvironmentVariable"
    "History Cleaner"            = "HistoryCleaner"
    "Registry Backup Restore"    = "RegistryBackupRestore"
    "Startup Application List"   = "StartupApplicationList"
    "Task Scheduler"             = "TaskScheduler"
    "Advanced Settings"          = "AdvancedSettings"
    "Change VNC Password"        = "ChangeVNCPassword"
    "General Settings"           = "GeneralSettings"
    "USB Device Manager"         = "USBDeviceManager"
}

foreach ($oldName in $renames.Keys) {
    $wb.Worksheets.Item($oldName).Name = $renames[$oldName]
}

# ---------------------------------------------------------------------------
# 2. Add the new "UserManagement" sheet right after "USBDeviceManager"
# ---------------------------------------------------------------------------
$usbSheet = $wb.Worksheets.Item("USBDeviceManager")
$userMgmt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $usbSheet)
$userMgmt.Name = "UserManagement"

# Column widths to roughly match the reference layout
$userMgmt.Columns.Item(1).ColumnWidth = 15.28515625
$userMgmt.Columns.Item(2).ColumnWidth = 10
$userMgmt.Columns.Item(3).ColumnWidth = 10
$userMgmt.Columns.Item(4).ColumnWidth = 11.42578125
$userMgmt.Columns.Item(5).ColumnWidth = 13.28515625
$userMgmt.Columns.Item(6).ColumnWidth = 19
$userMgmt.Columns.Item(7).ColumnWidth = 11.140625
$userMgmt.Columns.Item(8).ColumnWidth = 32.7109375
$userMgmt.Columns.Item(9).ColumnWidth = 22.5703125
$userMgmt.Columns.Item(10).ColumnWidth = 12

# Header row
$userMgmt.Range("A1").Value = "Template Name"
$userMgmt.Range("B1").Value = "Select Tab"
$userMgmt.Range("C1").Value = "Username"
$userMgmt.Range("D1").Value = "Password"
$userMgmt.Range("E1").Value = "Full Name"
$userMgmt.Range("F1").Value = "Description"
$userMgmt.Range("G1").Value = "Member Of"
$userMgmt.Range("H1").Value = "User Can Not Change The Password"
$userMgmt.Range("I1").Value = "Password Never Expires"
$userMgmt.Range("J1").Value = "Disable User"

# Row 2 - Add User (enabled restrictions)
$userMgmt.Range("A2").Value = "testingpr1"
$userMgmt.Range("B2").Value = "Add User"
$userMgmt.Range("C2").Value = "Admin"
$userMgmt.Hyperlinks.Add($userMgmt.Range("D2"), "http://www.test.com/", "", "", "admin@123")
$userMgmt.Range("E2").Value = "Admin Admin"
$userMgmt.Range("F2").Value = "testdata description"
$userMgmt.Range("G2").Value = "User"
$userMgmt.Range("H2").Value = "Y"
$userMgmt.Range("I2").Value = "Y"
$userMgmt.Range("J2").Value = "Y"

# Row 3 - Add User (disabled restrictions)
$userMgmt.Range("A3").Value = "testingpr1"
$userMgmt.Range("B3").Value = "Add User"
$userMgmt.Range("C3").Value = "Admin"
$userMgmt.Hyperlinks.Add($userMgmt.Range("D3"), "http://www.test.com/", "", "", "admin@123")
$userMgmt.Range("E3").Value = "Admin Admin"
$userMgmt.Range("F3").Value = "testdata description"
$userMgmt.Range("G3").Value = "User"
$userMgmt.Range("H3").Value = "N"
$userMgmt.Range("I3").Value = "N"
$userMgmt.Range("J3").Value = "N"

# Row 4 - Reset User
$userMgmt.Range("A4").Value = "testingpr1"
$userMgmt.Range("B4").Value = "Reset User"
$userMgmt.Range("C4").Value = "Admin"
$userMgmt.Hyperlinks.Add($userMgmt.Range("D4"), "http://www.test.com/", "", "", "admin@123")
$userMgmt.Range("E4").Value = "NA"
$userMgmt.Range("F4").Value = "NA"
$userMgmt.Range("G4").Value = "NA"
$userMgmt.Range("H4").Value = "NA"
$userMgmt.Range("I4").Value = "NA"
$userMgmt.Range("J4").Value = "NA"

# Shade the header row the same way the neighbouring sheets do (style 7 look alike)
$userMgmt.Range("A1:J1").Interior.ColorIndex = 6
$userMgmt.Range("A1:J1").Borders.LineStyle = 1

$userMgmt.Range("E1").Select()

# ---------------------------------------------------------------------------
# 3. Move "SoftwarePatchInstallUninstall" to the end, right before "Sheet1"
# ---------------------------------------------------------------------------
$patchSheet = $wb.Worksheets.Item("SoftwarePatchInstallUninstall")
$sheet1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$patchSheet.Move($sheet1)
$patchSheet.Select()
$patchSheet.Range("E15").Select()

foreach ($w in $wb.Worksheets) {
    Write-Output "$($w.Index): $($w.Name)"
}
